# Fruta / hortaliza, semanal
# Insert a new weekly price row above the current row 10 (the rest of the
# block's rows shift down by one), then populate it with this week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push existing rows 10-14 down to 11-15, leaving a fresh (empty) row 10.
$ws.Rows.Item(10).Insert()

# Fill in the new row 10 with the new weekly record.
$ws.Range("A10").Value = 7
$ws.Range("B10").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C10").Value = "Ñuble"
$ws.Range("D10").Value = 44868
$ws.Range("E10").Value = 16
$ws.Range("F10").Value = "Fruta"
$ws.Range("G10").Value = 100107
$ws.Range("H10").Value = "Otros"
$ws.Range("I10").Value = 100107002
$ws.Range("J10").Value = "Chirimoya"
$ws.Range("K10").Value = "Cultivar IV Región"
$ws.Range("L10").Value = "Especial"
$ws.Range("M10").Value = 60
$ws.Range("N10").Value = 26000
$ws.Range("O10").Value = 26000
$ws.Range("P10").Value = 26000
$ws.Range("Q10").Value = "$/bandeja 10 kilos"
$ws.Range("R10").Value = "Provincia de Limarí"
$ws.Range("S10").Value = 2600
$ws.Range("T10").Value = 10

Write-Output "done"
